# Change titles in second slide (session2 - Javascript & jQuery.pptx)
#
# Each of slides 3-11 carries its visible "title" text not in the
# (empty) p:ph type="title" placeholder, but in the body idx="1" (or
# idx="13") placeholder that is styled with the Impact font. This
# script rewrites those title runs to the new, more descriptive
# headings, and restores the ALL-CAPS emphasis on the jQuery selector
# slide.

$p = $ppt.ActivePresentation

# Slide 3: "Javascript" -> "What is Javascript"
$p.Slides.Item(3).Shapes.Item(2).TextFrame.TextRange.Text = "What is Javascript"

# Slide 4: "Javascript (continue)" -> "How to use javascript in html"
$p.Slides.Item(4).Shapes.Item(2).TextFrame.TextRange.Text = "How to use javascript in html"

# Slide 5: "Javascript" -> "function of javascript"
$p.Slides.Item(5).Shapes.Item(2).TextFrame.TextRange.Text = "function of javascript"

# Slide 6: "jQuery" -> "what is jQuery"
$p.Slides.Item(6).Shapes.Item(2).TextFrame.TextRange.Text = "what is jQuery"

# Slide 7: "jQuery" -> "how to use jQuery in html"
$p.Slides.Item(7).Shapes.Item(3).TextFrame.TextRange.Text = "how to use jQuery in html"

# Slide 8: "JQUERY" -> "how to use selector in JQUERY" (keep the all-caps emphasis)
$titleRange8 = $p.Slides.Item(8).Shapes.Item(2).TextFrame.TextRange
$titleRange8.Text = "how to use selector in JQUERY"
$titleRange8.Font.AllCaps = $true

# Slide 9: "jquery" -> "use jquery to operate html element"
$p.Slides.Item(9).Shapes.Item(2).TextFrame.TextRange.Text = "use jquery to operate html element"

# Slide 10: "jquery" -> "using ajax in jquery"
$p.Slides.Item(10).Shapes.Item(3).TextFrame.TextRange.Text = "using ajax in jquery"

# Slide 11: "jquery" -> "using ajax in jquery"
$p.Slides.Item(11).Shapes.Item(2).TextFrame.TextRange.Text = "using ajax in jquery"
